# Refresh the Price (D) / Volume(1h) (E) columns of the cryptos table
# with the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.325.56"
$ws.Range("E2").Value = "  -1.58%  "

$ws.Range("D3").Value = "3.507.58"
$ws.Range("E3").Value = "  -3.71%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.44"
$ws.Range("E5").Value = "  +2.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "553.87"
$ws.Range("E6").Value = "  -4.93%  "

$ws.Range("D7").Value = "3.500.88"
$ws.Range("E7").Value = "  -3.75%  "

$ws.Range("E8").Value = "  -2.06%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.656"
$ws.Range("E10").Value = "  -3.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "62.71"
$ws.Range("E11").Value = "  +10.90%  "

$ws.Range("E12").Value = "  -7.05%  "

$ws.Range("E13").Value = "  -7.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.85"
$ws.Range("E14").Value = "  -2.77%  "

$ws.Range("D15").Value = "4.067.57"
$ws.Range("E15").Value = "  -3.86%  "

$ws.Range("D16").Value = "3.507.92"
$ws.Range("E16").Value = "  -3.81%  "

$ws.Range("E17").Value = "  -1.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.46"
$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("D19").Value = "67.069.29"
$ws.Range("E19").Value = "  -1.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.84"
$ws.Range("E20").Value = "  -6.04%  "

$ws.Range("E21").Value = "  -5.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.04"
$ws.Range("E22").Value = "  -2.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.55"
$ws.Range("E23").Value = "  -4.27%  "

$ws.Range("E24").Value = "  -6.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.02"
$ws.Range("E25").Value = "  -3.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.94"
$ws.Range("E26").Value = "  +1.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.28"
$ws.Range("E27").Value = "  -3.10%  "

$ws.Range("E28").Value = "  -4.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.86"
$ws.Range("E29").Value = "  -3.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "686.16"
$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.05"
$ws.Range("E32").Value = "  -13.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.73"
$ws.Range("E33").Value = "  -4.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.77"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.112"
$ws.Range("E35").Value = "  -6.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.78"
$ws.Range("E36").Value = "  -9.42%  "

$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.399"
$ws.Range("E38").Value = "  -6.27%  "

$ws.Range("E39").Value = "  -4.88%  "

$ws.Range("D40").Value = "3.074.17"
$ws.Range("E40").Value = "  -4.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.00"
$ws.Range("E42").Value = "  -4.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  -9.06%  "

$ws.Range("E44").Value = "  -14.37%  "

$ws.Range("E45").Value = "  +5.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("E46").Value = "  -10.66%  "

$ws.Range("E47").Value = "  -5.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.128"
$ws.Range("E48").Value = "  -3.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.15"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.23"
$ws.Range("E50").Value = "  -7.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.88"
$ws.Range("E51").Value = "  -7.67%  "

